# test3.xlsx was re-saved by Excel; the only content-level change in the
# diff that the Excel object model can actually reproduce is the worksheet
# being renamed from "Tabelle1" to "Tabelle3". Everything else in the diff
# (the dropped mc:AlternateContent/x15ac:absPath block, the bookView window
# size, calcPr/calcId, the added <oleSize ref="A1"/>, and the
# sheetFormatPr/row x14ac:dyDescent shift from 0.3 to 0.25) are incidental
# artifacts of which Excel build/process resaved the file - there is no
# Workbook/Worksheet/Window property that drives any of that XML, so there
# is nothing meaningful to script for it.

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Tabelle1") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "Tabelle3"
